$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source added a new weekly price observation. This shifts the existing
# records at row 421 and below down by one row, and inserts the new record
# as row 421.
$ws.Rows.Item(421).EntireRow.Insert()

$ws.Cells.Item(421, 1).Value = 6
$ws.Cells.Item(421, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(421, 3).Value = "Metropolitana"
$ws.Cells.Item(421, 4).Value = 44610
$ws.Cells.Item(421, 5).Value = 13
$ws.Cells.Item(421, 6).Value = 100112044
$ws.Cells.Item(421, 7).Value = "Perejil"
$ws.Cells.Item(421, 8).Value = "Sin especificar"
$ws.Cells.Item(421, 9).Value = "Primera"
$ws.Cells.Item(421, 10).Value = 170
$ws.Cells.Item(421, 11).Value = 15000
$ws.Cells.Item(421, 12).Value = 16000
$ws.Cells.Item(421, 13).Value = 15353
$ws.Cells.Item(421, 14).Value = "`$/docena de atados"
$ws.Cells.Item(421, 15).Value = "Región Metropolitana"
$ws.Cells.Item(421, 16).Value = 5118
$ws.Cells.Item(421, 17).Value = 3
$ws.Cells.Item(421, 18).Value = "Hortaliza"
